# Refresh the "ランサーズ" (Lancers) job-listing sheet with a newer scrape
# (2025-09-29 06:27:26 JST), replacing the previous scrape's 10 data rows
# with this run's 6 data rows, and narrowing a few columns to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Drop every existing hyperlink (and its relationship) up front so
#        we don't end up with stale/duplicate relationship entries once the
#        row count shrinks and the URLs change. ---
$ws.Cells.Hyperlinks.Delete()

# --- 2. Overwrite rows 2-7 with the new scrape's data. ---
$timestamp = "2025-09-29 06:27:26"
$category = "システム開発"
$deadline = "期限情報なし"

$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "【急募】リスト抽出ツール開発のフリーランスを探しています!"
$ws.Range("C2").Value = $category
$ws.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E2").Value = $deadline
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5402362"
$ws.Range("G2").Value = 128
$ws.Range("H2").Value = "◆ツール,開発"

$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "日程を作成するシステムの開発"
$ws.Range("C3").Value = $category
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E3").Value = $deadline
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5402412"
$ws.Range("G3").Value = 78
$ws.Range("H3").Value = "◆開発"

$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "【ペットのアバター化】Pawsitiveプロトタイプ開発の依頼"
$ws.Range("C4").Value = $category
$ws.Range("D4").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E4").Value = $deadline
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5399313"
$ws.Range("G4").Value = 68
$ws.Range("H4").Value = "◆開発"

$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "【急募】新しい口コミサイトの構築をお手伝いください!"
$ws.Range("C5").Value = $category
$ws.Range("D5").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E5").Value = $deadline
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5402277"
$ws.Range("G5").Value = 38
$ws.Range("H5").Value = "◇サイト"

$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "限定公開 PR 限定公開の仕事"
$ws.Range("C6").Value = $category
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E6").Value = $deadline
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5399347"
$ws.Range("G6").Value = 13
$ws.Range("H6").ClearContents()

$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "エンジニア面談をお願い致します"
$ws.Range("C7").Value = $category
$ws.Range("D7").Value = "~ 5,000 円 / 固定"
$ws.Range("E7").Value = $deadline
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5402603"
$ws.Range("G7").Value = 10
$ws.Range("H7").ClearContents()

# --- 3. Drop the now-unused rows 8-11 (this also shrinks the sheet
#        dimension from A1:H11 down to A1:H7). ---
$ws.Range("A8:A11").EntireRow.Delete()

# --- 4. Re-create the hyperlinks for the surviving URL cells (F2:F7),
#        then reapply the built-in "Hyperlink" cell style so they match
#        the original look (same style index as before, no new style
#        entries minted). (Addresses are the same literals just written
#        into each cell above - avoids relying on a COM property getter.) ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5402362")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5402412")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5399313")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5402277")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5399347")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5402603")
$ws.Range("F2:F7").Style = "Hyperlink"

# --- 5. Narrow columns B, D and H (ColumnWidth is offset from the stored
#        <col width> by a fixed 5/6 padding, so subtract that to land on
#        the exact target widths of 33 / 28 / 12). ---
$ws.Columns.Item(2).ColumnWidth = 33 - 5/6
$ws.Columns.Item(4).ColumnWidth = 28 - 5/6
$ws.Columns.Item(8).ColumnWidth = 12 - 5/6
